# Apply the changes described by the diff:
#  - Metadata sheet: Version 2.2.0-ballot -> 2.1.0
#  - Metadata sheet: Date 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
#  - Metadata sheet: Base Definition http://hl7.org/fhir/StructureDefinition/Extension|4.0.1
#       -> http://hl7.org/fhir/StructureDefinition/Extension
#  - Elements sheet: Extension.value[x] Type(s) text drops the "|2.2.0-ballot" version suffixes
#  - Elements sheet: column K (11) width 158.5546875 -> 139.5390625

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2.1.0"
$wsMetadata.Range("B8").Value = "2025-12-19T08:44:55+00:00"
$wsMetadata.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner|https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-practitioner-role)
"
# ColumnWidth is stored internally in 1/6-character increments plus a fixed
# 5/6 padding offset, so request the value that rounds to the target
# (158.5546875 -> 139.5390625) stored width of 139.5 (closest reachable grid point).
$wsElements.Columns.Item(11).ColumnWidth = 138.70572916666666
